# This document was originally generated by a tool (e.g. pandoc) that wrote a
# short "Heading1" title paragraph followed by a bold "By <author>" paragraph.
# The new version instead uses pandoc-style title-block paragraph styles
# ("Title" and "Authors"), with the title/author text split into one run per
# word (and one run per inter-word space), and drops the leading "By " before
# the author's name.

$d = $word.ActiveDocument

# The first two paragraphs of the document are the ones being rewritten:
#   1) "Funds Needed To Carry On Work in N. Y." (currently styled Heading1)
#   2) "By Dorothy Day" (currently bold, unstyled)
$titlePara  = $d.Paragraphs.Item(1)
$authorPara = $d.Paragraphs.Item(2)

if ($titlePara.Range.Text.Trim() -ne "Funds Needed To Carry On Work in N. Y." -or `
    $authorPara.Range.Text.Trim() -ne "By Dorothy Day") {
    throw "Unexpected document content; expected the title/author paragraphs at positions 1 and 2."
}

# Range spanning both paragraphs (start of paragraph 1 through end of
# paragraph 2, i.e. including the paragraph mark that separates them but not
# a trailing one, so paragraph 3 onward is left completely untouched).
$targetRange = $d.Range($titlePara.Range.Start, $authorPara.Range.End)

# Build the replacement as a single run per "word" (pandoc style), matching
# the target markup exactly: plain runs with no rPr, separated by single
# space runs.
$titleWords  = @('Funds', ' ', 'Needed', ' ', 'To', ' ', 'Carry', ' ', 'On', ' ', 'Work', ' ', 'in', ' ', 'N', '.', ' ', 'Y', '.')
$authorWords = @('Dorothy', ' ', 'Day')

function Make-Runs($words) {
    $sb = ""
    foreach ($w in $words) {
        $sb += '<w:r><w:t xml:space="preserve">' + $w + '</w:t></w:r>'
    }
    return $sb
}

$titleRunsXml  = Make-Runs $titleWords
$authorRunsXml = Make-Runs $authorWords

$bodyXml = '<w:p><w:pPr><w:pStyle w:val="Title"/></w:pPr>' + $titleRunsXml + '</w:p>' `
         + '<w:p><w:pPr><w:pStyle w:val="Authors"/></w:pPr>' + $authorRunsXml + '</w:p>'

# InsertXML requires a full WordprocessingML package fragment so that runs
# are preserved as-authored instead of being merged back into one run the
# way plain Range.Text/InsertAfter edits would.
$xmlPackage = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' `
            + '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' `
            + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' `
            + '<pkg:xmlData>' `
            + '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' `
            + '<w:body>' + $bodyXml + '</w:body>' `
            + '</w:document>' `
            + '</pkg:xmlData></pkg:part></pkg:package>'

$targetRange.InsertXML($xmlPackage)

Write-Host "Title paragraph now:" $d.Paragraphs.Item(1).Style.NameLocal "->" $d.Paragraphs.Item(1).Range.Text
Write-Host "Author paragraph now:" $d.Paragraphs.Item(2).Style.NameLocal "->" $d.Paragraphs.Item(2).Range.Text
